$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Result")

# Insert a new row at row 8, shifting existing rows 8-27 down to 9-28
$ws.Rows("8:8").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Fill in the new row 8 content
$ws.Range("B8").Value = "Edit EBM Only Order"
$ws.Range("C8").Value = $false

$ws.Range("B8").Select()
